$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Summary sheet: update aggregate stats after trade #16 closes
# ------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.96
$summary.Range("B4").Value = -0.04
$summary.Range("B5").Value = -0.05
$summary.Range("B6").Value = 16
$summary.Range("B7").Value = 6
$summary.Range("B9").Value = 37.5

# ------------------------------------------------------------------
# Strategy Status sheet: update MarketMaking strategy row (row 4)
# ------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.95999999999999
$status.Range("D4").Value = 16
$status.Range("E4").Value = -0.04
$status.Range("F4").Value = -0.04
$status.Range("G4").Value = 37.5

# ------------------------------------------------------------------
# All Trades + MarketMaking sheets: append closed trade #16 as row 17
# ------------------------------------------------------------------
$tradeSheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $tradeSheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Force text storage for date/time-looking strings so Excel does not
    # reinterpret them as date/time serial numbers, then restore the
    # default "Normal" style so no stray number format lingers on the cell.
    $ws.Range("B17").NumberFormat = "@"
    $ws.Range("C17").NumberFormat = "@"

    $ws.Range("A17").Value = 16
    $ws.Range("B17").Value = "2026-02-17"
    $ws.Range("C17").Value = "07:54:14"
    $ws.Range("B17").Style = "Normal"
    $ws.Range("C17").Style = "Normal"
    $ws.Range("D17").Value = "MarketMaking"
    $ws.Range("E17").Value = "UP"
    $ws.Range("F17").Value = 0.9399999999999999
    $ws.Range("G17").Value = 0.97
    $ws.Range("H17").Value = "CLOSED"
    $ws.Range("I17").Value = 3.1915
    $ws.Range("J17").Value = 0.03
    $ws.Range("K17").Value = 99.95999999999999
    $ws.Range("L17").Value = 0
    $ws.Range("M17").Value = 0
    $ws.Range("N17").Value = 0.6
    $ws.Range("O17").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P17").Value = "early_exit"
    $ws.Range("Q17").Value = 0.11
}
